$d = $word.ActiveDocument

# The first table is the work-log table (Date / Activity / Time columns).
$t = $d.Tables.Item(1)

# Row: 25/3/24 - troubleshooting / reimporting raspberry pi image / OpenPLC & Modbus testing - 3 hours
$row1 = $t.Rows.Add()
$row1.Cells.Item(1).Range.Text = "25/3/24"
$row1.Cells.Item(2).Range.Text = "Spent time troubleshooting. Had to reimport new image of raspberry pi clone into virtualbox due to OpenPLC connectivity issues. This resolved the earlier problem of the circuit not working. Created new circuit in OpenPLC editor and imported to OpenPLC, tested with Modbus with success."
$row1.Cells.Item(3).Range.Text = "3"

# Row: 25/3/24 - Node-Red flow work / OpenPLC edits - 4 hours
$row2 = $t.Rows.Add()
$row2.Cells.Item(1).Range.Text = "25/3/24"
$row2.Cells.Item(2).Range.Text = "Worked on Node-Red flow to OpenPLC to check Modbus reading and UI response. Made edits to OpenPLC to resolve issues"
$row2.Cells.Item(3).Range.Text = "4"
